$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "この古くから伝わるレシピは、一口飲むごとに本格的でしっかりとした味わいを約束します。"; New = "Esta receita centenária promete um sabor autêntico e robusto em cada gole." },
    @{ Old = "風味は強烈でありながらバランスが取れており、快適で心地よい体験を生み出します。"; New = "Os sabores são intensos, mas equilibrados, proporcionando uma experiência reconfortante e relaxante." },
    @{ Old = "お好みの方法でチャイをお楽しみいただけるよう、簡単な淹れ方の説明書が付属しています。"; New = "Instruções simples de preparo estão incluídas para ajudá-lo a saborear seu chai exatamente do jeito que você gosta." },
    @{ Old = "Mystic Spice Chai Tea がお客様のご期待に添えない場合は、当社が改善するよう努めます。"; New = "Se o Mystic Spice Chai Tea não atender suas expectativas, estamos comprometidos em resolver da melhor maneira possível." }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
